# feat: add 2022-Q4 data
#
# The workbook has three sheets: "总计" (totals), "2022-Q3", "2021-Q2".
# This adds a new "2022-Q4" sheet (with its fund-holdings detail, copied
# in shape from the existing "2022-Q3" sheet) positioned right after
# "总计" and before "2022-Q3", and inserts a corresponding summary row
# into "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row for 2022-Q4 above the 2022-Q3 row,
#    pushing the existing 2022-Q3 / 2021-Q2 rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift row 3 (2021-Q2) down to row 4, then row 2 (2022-Q3) down to row 3.
# Copying whole rows (not just setting .Value) keeps the existing cell
# styles (e.g. the centred/bordered style on column A) intact.
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

# New first data row: 2022-Q4 summary.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 1.19

# Fix up the running index column now that the rows moved.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 2) Add the "2022-Q4" worksheet itself, placed before "2022-Q3".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Seed the new sheet from "2022-Q3" so column layout/styles match
# (header row style + the bordered/centred style used on column A).
$q3.Range("A1:H3").Copy($q4.Range("A1:H3"))

# Fund-size / percentage columns are stored as plain text in this
# workbook (matches the convention used on the other quarter sheets),
# so force text formatting before writing the numeric-looking strings.
$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "012526"
$q4.Range("C2").Value = "广发盛锦混合A"
$q4.Range("D2").Value = "24.00"
$q4.Range("E2").Value = "93.13"
$q4.Range("F2").Value = "4.74"
$q4.Range("G2").Value = "1.1376"
$q4.Range("H2").Value = 6

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "012527"
$q4.Range("C3").Value = "广发盛锦混合C"
$q4.Range("D3").Value = "1.14"
$q4.Range("E3").Value = "93.13"
$q4.Range("F3").Value = "4.74"
$q4.Range("G3").Value = "0.0540"
$q4.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 3) Keep "2021-Q2" (now the last sheet) the active/selected tab, as it
#    was before the new sheet was inserted.
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
